# Update mods data [2026-01-17 15:09:23]
# Append a new row (68) to the ModCounts sheet with the next day's data,
# mirroring the formatting of the previous row (67).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing row onto the new row
# so the new cells share the same style index as the rest of the table.
$ws.Range("A67:C67").Copy()
$ws.Range("A68:C68").PasteSpecial(-4122)

# Column A holds a literal date-like string (not a real Excel date).
# Setting the NumberFormat to Text ("@") before assigning the value stops
# Excel from auto-converting the "2026/01/17" string into a date serial
# number. Afterwards, restore the "Normal" style plus the original
# centered alignment so the cell's style matches the rest of the column.
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "2026/01/17"
$ws.Range("A68").Style = "Normal"
$ws.Range("A68").HorizontalAlignment = -4108
$ws.Range("A68").VerticalAlignment = -4108

$ws.Range("B68").Value = "逃离鸭科夫"
$ws.Range("C68").Value = 1147
